$wb = $excel.ActiveWorkbook

$urlMd09709514 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/429c14039e2ed7bf7b7f95a19f7eb8a1462d079b/e2e/09709514-9084-498d-8151-da8e1bc70d92.md"
$urlMdA6f5cf38 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/429c14039e2ed7bf7b7f95a19f7eb8a1462d079b/e2e/a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: handback status for both languages is now in sync
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$wsOverview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: record that the handback for both files has happened
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Range("I2").Value = "09709514-9084-498d-8151-da8e1bc70d92.md"
$wsZhCn.Range("J2").Value = "09709514-9084-498d-8151-da8e1bc70d92.7c4fde9576d957805e17575fd07c357f46436879.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-22 16:48:44"

$wsZhCn.Range("I3").Value = "a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md"
$wsZhCn.Range("J3").Value = "a6f5cf38-7a53-43f0-9081-cb5230a1fab6.4d705ba9d7e725c8bf5e2e4c7bb802ccb8105142.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-22 16:48:44"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlMd09709514, "", "", "09709514-9084-498d-8151-da8e1bc70d92.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlMdA6f5cf38, "", "", "a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md") | Out-Null
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Range("I3").Style = "HyperLink"

$wsZhCn.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$wsZhCn.Columns.Item(9).EntireColumn.AutoFit() | Out-Null
$wsZhCn.Columns.Item(10).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: record that the handback for both files has happened
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Range("I2").Value = "09709514-9084-498d-8151-da8e1bc70d92.md"
$wsDeDe.Range("J2").Value = "09709514-9084-498d-8151-da8e1bc70d92.7c4fde9576d957805e17575fd07c357f46436879.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-22 16:48:50"

$wsDeDe.Range("I3").Value = "a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md"
$wsDeDe.Range("J3").Value = "a6f5cf38-7a53-43f0-9081-cb5230a1fab6.4d705ba9d7e725c8bf5e2e4c7bb802ccb8105142.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-22 16:48:50"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlMd09709514, "", "", "09709514-9084-498d-8151-da8e1bc70d92.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlMdA6f5cf38, "", "", "a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md") | Out-Null
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Range("I3").Style = "HyperLink"

$wsDeDe.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$wsDeDe.Columns.Item(9).EntireColumn.AutoFit() | Out-Null
$wsDeDe.Columns.Item(10).EntireColumn.AutoFit() | Out-Null
